$wb = $excel.ActiveWorkbook

# --- 1. Create the new "types" lookup sheet and position it before "Sheet1" ---
$sheet1 = $wb.Worksheets.Item("Sheet1")
$typesWs = $wb.Worksheets.Add()
$typesWs.Name = "types"
$typesWs.Move($sheet1)

# --- 2. Populate the "types" sheet with the lookup table ---
$typesData = @(
    @(1, "plane", "avion"),
    @(2, "boat", "bateau"),
    @(3, "drone", "drone"),
    @(4, "underwater glider", "planeur sous-marin"),
    @(5, "land", "terre"),
    @(6, "mooring", "mouillage"),
    @(7, "space", "espace")
)

for ($i = 0; $i -lt $typesData.Length; $i++) {
    $row = $typesData[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $typesWs.Cells.Item($i + 1, $j + 1).Value = $row[$j]
    }
}

$typesWs.Range("A1:C7").Select() | Out-Null

# --- 3. Update the "Sheet1" data sheet ---
$sheet1 = $wb.Worksheets.Item("Sheet1")

# header for new column F
$sheet1.Cells.Item(1, 6).Value = "longname"

# row 13 type changes from "hydrophone" to "mooring"
$sheet1.Cells.Item(13, 1).Value = "mooring"

# --- 4. Make "Sheet1" the active sheet/tab, with selection A21 ---
$sheet1.Activate() | Out-Null
$sheet1.Range("A21").Select() | Out-Null
